$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the "Zapallo italiano" price series.
# It needs to land as the new row 390 (in date order within the sheet),
# which pushes the former rows 390:471 down to 391:472.
$ws.Rows.Item(390).Insert()

# Populate the newly inserted row 390 with the new record's values
$ws.Cells.Item(390, 1).Value = 4
$ws.Cells.Item(390, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(390, 3).Value = 'Los Lagos'
$ws.Cells.Item(390, 4).Value = 45204
$ws.Cells.Item(390, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(390, 5).Value = 10
$ws.Cells.Item(390, 6).Value = 100112032
$ws.Cells.Item(390, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(390, 8).Value = 'Sin especificar'
$ws.Cells.Item(390, 9).Value = 'Primera'
$ws.Cells.Item(390, 10).Value = 80
$ws.Cells.Item(390, 11).Value = 26000
$ws.Cells.Item(390, 12).Value = 26000
$ws.Cells.Item(390, 13).Value = 26000
$ws.Cells.Item(390, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(390, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(390, 16).Value = 520
$ws.Cells.Item(390, 17).Value = 50
$ws.Cells.Item(390, 18).Value = 'Hortaliza'
